$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.424.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "'2.524.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'573.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'166.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "'2.524.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("D11").Value = "'0.167"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  +4.39%  "
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").Value = "'2.987.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "'69.185.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "'24.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "'2.531.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "'11.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'7.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").Value = "'349.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'2.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'70.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("E26").Value = "  -2.22%  "
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").Value = "'7.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").Value = "'464.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.99%  "
$ws.Range("E33").Value = "  -3.70%  "
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").Value = "'157.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "'19.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").Value = "'18.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "'38.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  -4.17%  "
$ws.Range("E46").Value = "  -12.75%  "
$ws.Range("D47").Value = "'142.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -2.52%  "
